# Auto-generated Excel COM-interop script applying the Maduin_Profits.xlsx diff.
# For each affected leve-profit row, update the recomputed market-price /
# profit columns (H:N) per sheet to match the latest price-scan run.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 319.85715
$ws.Range("J2").Value = 725
$ws.Range("L2").Value = 725
$ws.Range("N2").Value = -951
# Row 43
$ws.Range("H43").Value = 8357
$ws.Range("I43").Value = 9666.667
$ws.Range("J43").Value = 7374.75
$ws.Range("K43").Value = 9666.667
$ws.Range("L43").Value = 7374.75
$ws.Range("M43").Value = -9597.667
$ws.Range("N43").Value = -7512.75
# Row 62
$ws.Range("H62").Value = 4891.8887
$ws.Range("I62").Value = 5518.8335
$ws.Range("K62").Value = 5518.8335
$ws.Range("M62").Value = -4894.8335
# Row 65
$ws.Range("H65").Value = 4891.8887
$ws.Range("I65").Value = 5518.8335
$ws.Range("K65").Value = 27594.1675
$ws.Range("M65").Value = -24474.1675
# Row 86
$ws.Range("H86").Value = 6867.3335
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 6867.3335
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 6867.3335
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -9113.3335
# Row 89
$ws.Range("H89").Value = 6867.3335
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 6867.3335
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 34336.6675
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -45568.6675

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 3
$ws.Range("H3").Value = 2649
$ws.Range("I3").Value = 148.5
$ws.Range("J3").Value = 7650
$ws.Range("K3").Value = 148.5
$ws.Range("L3").Value = 7650
$ws.Range("M3").Value = -33.5
$ws.Range("N3").Value = -7880
# Row 8
$ws.Range("H8").Value = 4012651
$ws.Range("I8").Value = 6683601.5
$ws.Range("J8").Value = 6225
$ws.Range("K8").Value = 6683601.5
$ws.Range("L8").Value = 6225
$ws.Range("M8").Value = -6683457.5
$ws.Range("N8").Value = -6513
# Row 17
$ws.Range("H17").Value = 5002.6665
$ws.Range("I17").Value = 5002.6665
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 5002.6665
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -4829.6665
$ws.Range("N17").ClearContents()
# Row 24
$ws.Range("H24").Value = 61379.145
$ws.Range("J24").Value = 61379.145
$ws.Range("L24").Value = 61379.145
$ws.Range("N24").Value = -62127.145
# Row 95
$ws.Range("H95").Value = 44666.668
$ws.Range("I95").Value = 49000
$ws.Range("J95").Value = 43800
$ws.Range("K95").Value = 49000
$ws.Range("L95").Value = 43800
$ws.Range("M95").Value = -46254
$ws.Range("N95").Value = -49292
# Row 100
$ws.Range("H100").Value = 61379.145
$ws.Range("J100").Value = 61379.145
$ws.Range("L100").Value = 61379.145
$ws.Range("N100").Value = -63543.145

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 3167028
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
# Row 10
$ws.Range("H10").Value = 2999
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 2999
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 2999
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -3279
# Row 12
$ws.Range("H12").Value = 1046.1428
$ws.Range("I12").Value = 724.6667
$ws.Range("J12").Value = 1287.25
$ws.Range("K12").Value = 724.6667
$ws.Range("L12").Value = 1287.25
$ws.Range("M12").Value = -556.6667
$ws.Range("N12").Value = -1623.25
# Row 14
$ws.Range("H14").Value = 7999.9
$ws.Range("I14").Value = 9700
$ws.Range("K14").Value = 9700
$ws.Range("M14").Value = -9528
# Row 16
$ws.Range("H16").Value = 274.5
$ws.Range("I16").Value = 166
$ws.Range("J16").Value = 600
$ws.Range("K16").Value = 166
$ws.Range("L16").Value = 600
$ws.Range("M16").Value = 4
$ws.Range("N16").Value = -940
# Row 18
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
# Row 86
$ws.Range("H86").Value = 1518.1
$ws.Range("I86").Value = 1424.5
$ws.Range("J86").Value = 1892.5
$ws.Range("K86").Value = 1424.5
$ws.Range("L86").Value = 1892.5
$ws.Range("M86").Value = -301.5
$ws.Range("N86").Value = -4138.5
# Row 89
$ws.Range("H89").Value = 1518.1
$ws.Range("I89").Value = 1424.5
$ws.Range("J89").Value = 1892.5
$ws.Range("K89").Value = 7122.5
$ws.Range("L89").Value = 9462.5
$ws.Range("M89").Value = -1506.5
$ws.Range("N89").Value = -20694.5
# Row 92
$ws.Range("H92").Value = 20000
$ws.Range("J92").Value = 20000
$ws.Range("L92").Value = 20000
$ws.Range("N92").Value = -24992

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 6410.5
$ws.Range("I4").Value = 1821.4
$ws.Range("J4").Value = 10999.6
$ws.Range("K4").Value = 1821.4
$ws.Range("L4").Value = 10999.6
$ws.Range("M4").Value = -1709.4
$ws.Range("N4").Value = -11223.6
# Row 88
$ws.Range("H88").Value = 31780.666
$ws.Range("J88").Value = 31780.666
$ws.Range("L88").Value = 31780.666
$ws.Range("N88").Value = -32592.666
# Row 91
$ws.Range("H91").Value = 31780.666
$ws.Range("J91").Value = 31780.666
$ws.Range("L91").Value = 31780.666
$ws.Range("N91").Value = -34588.666

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 250003740
$ws.Range("J4").Value = 5000
$ws.Range("L4").Value = 15000
$ws.Range("N4").Value = -15224
# Row 122
$ws.Range("H122").Value = 1150.909
$ws.Range("I122").Value = 624.6667
$ws.Range("K122").Value = 5622.0003
$ws.Range("M122").Value = -3172.0003

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 500
$ws.Range("I5").Value = 500
$ws.Range("K5").Value = 500
$ws.Range("M5").Value = -388
# Row 6
$ws.Range("H6").Value = 1700.75
$ws.Range("J6").Value = 1901.5
$ws.Range("L6").Value = 1901.5
$ws.Range("N6").Value = -2127.5
# Row 12
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
# Row 13
$ws.Range("H13").Value = 131.6
$ws.Range("I13").Value = 131.6
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 131.6
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 7.400000000000006
$ws.Range("N13").ClearContents()
# Row 16
$ws.Range("H16").Value = 1700.75
$ws.Range("J16").Value = 1901.5
$ws.Range("L16").Value = 1901.5
$ws.Range("N16").Value = -2401.5
# Row 122
$ws.Range("H122").Value = 1986.1428
$ws.Range("I122").Value = 1986.1428
$ws.Range("K122").Value = 5958.428400000001
$ws.Range("M122").Value = -3508.428400000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 5500
$ws.Range("I68").Value = 3000
$ws.Range("K68").Value = 3000
$ws.Range("M68").Value = -2251
# Row 71
$ws.Range("H71").Value = 5500
$ws.Range("I71").Value = 3000
$ws.Range("K71").Value = 15000
$ws.Range("M71").Value = -11256

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
# Row 12
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()

